$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fb_senerios")

# Rename locator "id=login" to "name=login"
$ws.Range("B6").Value = "name=login"

# Fix the "enter   url" (multiple spaces) value to "enter url" (single space)
$ws.Range("C3").Value = "enter url"

# Update the active cell selection to C3 as per the edit
$ws.Range("C3").Select()
